# Regenerate save_data column G ("K" - strikeouts) to use actual K counts
# instead of the previous Strike# totals, for gallen_zac 2021 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values for rows 2-28, replacing the old Strike# values.
$newK = @{
    2  = 6
    3  = 8
    4  = 7
    5  = 4
    6  = 8
    7  = 4
    8  = 6
    9  = 9
    10 = 6
    11 = 9
    12 = 6
    13 = 8
    14 = 7
    15 = 6
    16 = 3
    17 = 4
    18 = 8
    19 = 3
    20 = 5
    21 = 7
    22 = 6
    23 = 6
    24 = 9
    25 = 3
    26 = 4
    27 = 3
    28 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
